# Generate Report for Handback
# - Update the "Ready for handoff" status text (shared across sheets) to
#   "Handback transform failed"
# - Add an Error Detail (column L) entry on the zh-cn and de-de sheets
#   describing the handback/handoff file-name mismatch.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Status text change - this string is shared across every sheet that
# references it (Overview B3/C3, zh-cn C3, de-de C3), so update each cell.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zh.Range("C3").Value = "Handback transform failed"
$de.Range("C3").Value = "Handback transform failed"

# New "Error Detail" values in column L, row 3, for each language sheet.
$zh.Range("L3").Value = "Handback file name: 1vknwv3z.zhu is different with handoff file name: c1a72339-8f1d-4830-a1cc-e552e69fa60b.780df331bf46f399b93b89df1a55116badb01f90.zh-cn."
$de.Range("L3").Value = "Handback file name: 1vknwv3z.zhu is different with handoff file name: c1a72339-8f1d-4830-a1cc-e552e69fa60b.780df331bf46f399b93b89df1a55116badb01f90.de-de."
